# Update crypto price list (D = Price, E = Volume(1h)) per latest scrape.
# D-column values are free-form numeric-looking TEXT (dotted thousands, many
# decimals, subscript notation, etc.) in the source data, so each write forces
# a text number-format before assignment and then clears formatting again (so
# the cell keeps its original default style) -- this avoids Excel's automatic
# "numeric-looking string -> number" coercion from corrupting values such as
# "571.47" (-> float) or "0.0000174" (-> scientific notation).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.410.81"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -1.24%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.432.83"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.64%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.27%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "571.47"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.09%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.34"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.52%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("E8").Value = "  -1.14%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.427.87"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.97%  "

$ws.Range("E10").Value = "  -4.07%  "

$ws.Range("E11").Value = "  +0.90%  "

$ws.Range("E12").Value = "  -1.04%  "

$ws.Range("E13").Value = "  -1.70%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.52"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.71%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000174"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -3.91%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.144.71"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -1.69%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.410.36"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.39%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.05"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -2.70%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.12"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.75%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "325.31"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.46%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.13"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.72%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.02"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +4.01%  "

$ws.Range("E24").Value = "  +0.19%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.25"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -2.16%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "621.80"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.13%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.98"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.46%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0963"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -6.88%  "

$ws.Range("E30").Value = "  +0.30%  "

$ws.Range("E31").Value = "  -3.44%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.02"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -3.33%  "

$ws.Range("E33").Value = "  -2.01%  "

$ws.Range("E34").Value = "  -7.87%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.01"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -2.16%  "

$ws.Range("E36").Value = "  +0.13%  "

$ws.Range("E37").Value = "  -4.51%  "

$ws.Range("E38").Value = "  -2.00%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.44"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -1.75%  "

$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.24"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -2.87%  "

$ws.Range("B41").Value = "Monero"
$ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "146.17"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.25%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.73"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -5.04%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "42.31"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +1.01%  "

$ws.Range("E44").Value = "  -0.02%  "

$ws.Range("E45").Value = "  -6.37%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "145.28"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -2.23%  "

$ws.Range("E47").Value = "  -1.40%  "

$ws.Range("E48").Value = "  -3.81%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.594"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.25%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.81"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -4.82%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0229"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -3.09%  "
